# Fix the "客単価" (H column) calculation on the "ABC分析_客構成" sheet.
# Previously H was computed as B / C; it should be B / (C * E).
# Row 10 is left untouched because its H value does not follow this
# relationship (the original data already had a different, explicit value).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ABC分析_客構成")

$rows = @(2,3,4,5,6,7,8,9,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,36,37,38,39,40,41,42,43,44,45,46,47,48,49,50,51)

foreach ($r in $rows) {
    $b = $ws.Cells.Item($r, 2).Value()
    $c = $ws.Cells.Item($r, 3).Value()
    $e = $ws.Cells.Item($r, 5).Value()
    $ws.Cells.Item($r, 8).Value = $b / ($c * $e)
}
